$wb = $excel.ActiveWorkbook

# Add the two new example sheets, positioned right after "Joe".
$joe = $wb.Worksheets.Item(1)
$debts = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $joe)
$debts.Name = "Debts"
$fixedAssets = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $debts)
$fixedAssets.Name = "Fixed Assets"

# --- Debts sheet headers ---
$debts.Range("A1").Value = "name"
$debts.Range("B1").Value = "type"
$debts.Range("C1").Value = "year"
$debts.Range("D1").Value = "term"
$debts.Range("E1").Value = "amount"
$debts.Range("F1").Value = "rate"
$debts.Range("A1:F1").Font.Bold = $true
$debts.Range("E1").NumberFormat = """$""#,##0"
[void]$debts.Range("A1:XFD1").Select()

# --- Fixed Assets sheet headers ---
$fixedAssets.Range("A1").Value = "name"
$fixedAssets.Range("B1").Value = "type"
$fixedAssets.Range("C1").Value = "basis"
$fixedAssets.Range("D1").Value = "value"
$fixedAssets.Range("E1").Value = "rate"
$fixedAssets.Range("F1").Value = "yod"
$fixedAssets.Range("G1").Value = "commission"
$fixedAssets.Range("A1:G1").Font.Bold = $true
$fixedAssets.Range("C1:D1").NumberFormat = """$""#,##0"
$fixedAssets.Columns.Item(7).ColumnWidth = 11.5
[void]$fixedAssets.Range("I6").Select()

# "Fixed Assets" is the sheet that was active/visible when the workbook was
# last saved.
$fixedAssets.Activate()
